# Update NATMI LR-pair statistics for Tgfb3-Tgfbr1 (per Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 1.987640333333333
$ws.Range("H2").Value2 = 5.962921
$ws.Range("I2").Value2 = 0.0278174819837782
$ws.Range("J2").Value2 = 0.0278174819837782
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 11.19524
$ws.Range("N2").Value2 = 33.58572
$ws.Range("O2").Value2 = 0.07847249539938134
$ws.Range("P2").Value2 = 0.07847249539938135
$ws.Range("Q2").Value2 = 22.25211056534667
$ws.Range("R2").Value2 = 200.26899508812
$ws.Range("S2").Value2 = 0.002182907226994408
$ws.Range("T2").Value2 = 0.002182907226994409

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 1.987640333333333
$ws.Range("H3").Value2 = 5.962921
$ws.Range("I3").Value2 = 0.0278174819837782
$ws.Range("J3").Value2 = 0.0278174819837782
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 31.07813833333333
$ws.Range("N3").Value2 = 93.234415
$ws.Range("O3").Value2 = 0.2178407133195748
$ws.Range("P3").Value2 = 0.2178407133195749
$ws.Range("Q3").Value2 = 61.77216123624611
$ws.Range("R3").Value2 = 555.949451126215
$ws.Range("S3").Value2 = 0.006059780118100665
$ws.Range("T3").Value2 = 0.006059780118100666

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 1.987640333333333
$ws.Range("H4").Value2 = 5.962921
$ws.Range("I4").Value2 = 0.0278174819837782
$ws.Range("J4").Value2 = 0.0278174819837782
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 91.56894199999999
$ws.Range("N4").Value2 = 274.706826
$ws.Range("O4").Value2 = 0.6418480872068143
$ws.Range("P4").Value2 = 0.6418480872068144
$ws.Range("Q4").Value2 = 182.0061223998606
$ws.Range("R4").Value2 = 1638.055101598746
$ws.Range("S4").Value2 = 0.01785459760219806
$ws.Range("T4").Value2 = 0.01785459760219806

# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 1.987640333333333
$ws.Range("H5").Value2 = 5.962921
$ws.Range("I5").Value2 = 0.0278174819837782
$ws.Range("J5").Value2 = 0.0278174819837782
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 8.822188333333335
$ws.Range("N5").Value2 = 26.466565
$ws.Range("O5").Value2 = 0.06183870407422939
$ws.Range("P5").Value2 = 0.06183870407422939
$ws.Range("Q5").Value2 = 17.53533735959611
$ws.Range("R5").Value2 = 157.818036236365
$ws.Range("S5").Value2 = 0.001720197036485068
$ws.Range("T5").Value2 = 0.001720197036485068

# Row 6
$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 56.81334033333334
$ws.Range("H6").Value2 = 170.440021
$ws.Range("I6").Value2 = 0.7951157181995667
$ws.Range("J6").Value2 = 0.7951157181995667
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 11.19524
$ws.Range("N6").Value2 = 33.58572
$ws.Range("O6").Value2 = 0.07847249539938134
$ws.Range("P6").Value2 = 0.07847249539938135
$ws.Range("Q6").Value2 = 636.0389802333467
$ws.Range("R6").Value2 = 5724.35082210012
$ws.Range("S6").Value2 = 0.06239471453839129
$ws.Range("T6").Value2 = 0.0623947145383913

# Row 7
$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 56.81334033333334
$ws.Range("H7").Value2 = 170.440021
$ws.Range("I7").Value2 = 0.7951157181995667
$ws.Range("J7").Value2 = 0.7951157181995667
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 31.07813833333333
$ws.Range("N7").Value2 = 93.234415
$ws.Range("O7").Value2 = 0.2178407133195748
$ws.Range("P7").Value2 = 0.2178407133195749
$ws.Range("Q7").Value2 = 1765.652850058079
$ws.Range("R7").Value2 = 15890.87565052272
$ws.Range("S7").Value2 = 0.1732085752241996
$ws.Range("T7").Value2 = 0.1732085752241997

# Row 8
$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 56.81334033333334
$ws.Range("H8").Value2 = 170.440021
$ws.Range("I8").Value2 = 0.7951157181995667
$ws.Range("J8").Value2 = 0.7951157181995667
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 91.56894199999999
$ws.Range("N8").Value2 = 274.706826
$ws.Range("O8").Value2 = 0.6418480872068143
$ws.Range("P8").Value2 = 0.6418480872068144
$ws.Range("Q8").Value2 = 5202.337465809261
$ws.Range("R8").Value2 = 46821.03719228334
$ws.Range("S8").Value2 = 0.5103435028344643
$ws.Range("T8").Value2 = 0.5103435028344644

# Row 9
$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 56.81334033333334
$ws.Range("H9").Value2 = 170.440021
$ws.Range("I9").Value2 = 0.7951157181995667
$ws.Range("J9").Value2 = 0.7951157181995667
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 8.822188333333335
$ws.Range("N9").Value2 = 26.466565
$ws.Range("O9").Value2 = 0.06183870407422939
$ws.Range("P9").Value2 = 0.06183870407422939
$ws.Range("Q9").Value2 = 501.2179882664295
$ws.Range("R9").Value2 = 4510.961894397866
$ws.Range("S9").Value2 = 0.04916892560251137
$ws.Range("T9").Value2 = 0.04916892560251138

# Row 10
$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 0.7501196666666666
$ws.Range("H10").Value2 = 2.250359
$ws.Range("I10").Value2 = 0.01049809664416703
$ws.Range("J10").Value2 = 0.01049809664416703
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 11.19524
$ws.Range("N10").Value2 = 33.58572
$ws.Range("O10").Value2 = 0.07847249539938134
$ws.Range("P10").Value2 = 0.07847249539938135
$ws.Range("Q10").Value2 = 8.397769697053333
$ws.Range("R10").Value2 = 75.57992727348001
$ws.Range("S10").Value2 = 0.0008238118406116581
$ws.Range("T10").Value2 = 0.0008238118406116583

# Row 11
$ws.Range("E11").Value2 = 3
$ws.Range("G11").Value2 = 0.7501196666666666
$ws.Range("H11").Value2 = 2.250359
$ws.Range("I11").Value2 = 0.01049809664416703
$ws.Range("J11").Value2 = 0.01049809664416703
$ws.Range("K11").Value2 = 3
$ws.Range("M11").Value2 = 31.07813833333333
$ws.Range("N11").Value2 = 93.234415
$ws.Range("O11").Value2 = 0.2178407133195748
$ws.Range("P11").Value2 = 0.2178407133195749
$ws.Range("Q11").Value2 = 23.31232276722055
$ws.Range("R11").Value2 = 209.810904904985
$ws.Range("S11").Value2 = 0.002286912861463181
$ws.Range("T11").Value2 = 0.002286912861463182

# Row 12
$ws.Range("E12").Value2 = 3
$ws.Range("G12").Value2 = 0.7501196666666666
$ws.Range("H12").Value2 = 2.250359
$ws.Range("I12").Value2 = 0.01049809664416703
$ws.Range("J12").Value2 = 0.01049809664416703
$ws.Range("K12").Value2 = 3
$ws.Range("M12").Value2 = 91.56894199999999
$ws.Range("N12").Value2 = 274.706826
$ws.Range("O12").Value2 = 0.6418480872068143
$ws.Range("P12").Value2 = 0.6418480872068144
$ws.Range("Q12").Value2 = 68.68766425005933
$ws.Range("R12").Value2 = 618.188978250534
$ws.Range("S12").Value2 = 0.006738183250370886
$ws.Range("T12").Value2 = 0.006738183250370888

# Row 13
$ws.Range("E13").Value2 = 3
$ws.Range("G13").Value2 = 0.7501196666666666
$ws.Range("H13").Value2 = 2.250359
$ws.Range("I13").Value2 = 0.01049809664416703
$ws.Range("J13").Value2 = 0.01049809664416703
$ws.Range("K13").Value2 = 3
$ws.Range("M13").Value2 = 8.822188333333335
$ws.Range("N13").Value2 = 26.466565
$ws.Range("O13").Value2 = 0.06183870407422939
$ws.Range("P13").Value2 = 0.06183870407422939
$ws.Range("Q13").Value2 = 6.617696971870556
$ws.Range("R13").Value2 = 59.55927274683501
$ws.Range("S13").Value2 = 0.0006491886917213058
$ws.Range("T13").Value2 = 0.0006491886917213059

# Row 14
$ws.Range("E14").Value2 = 3
$ws.Range("G14").Value2 = 11.90182033333333
$ws.Range("H14").Value2 = 35.705461
$ws.Range("I14").Value2 = 0.166568703172488
$ws.Range("J14").Value2 = 0.166568703172488
$ws.Range("K14").Value2 = 3
$ws.Range("M14").Value2 = 11.19524
$ws.Range("N14").Value2 = 33.58572
$ws.Range("O14").Value2 = 0.07847249539938134
$ws.Range("P14").Value2 = 0.07847249539938135
$ws.Range("Q14").Value2 = 133.2437350685467
$ws.Range("R14").Value2 = 1199.19361561692
$ws.Range("S14").Value2 = 0.01307106179338398
$ws.Range("T14").Value2 = 0.01307106179338398

# Row 15
$ws.Range("E15").Value2 = 3
$ws.Range("G15").Value2 = 11.90182033333333
$ws.Range("H15").Value2 = 35.705461
$ws.Range("I15").Value2 = 0.166568703172488
$ws.Range("J15").Value2 = 0.166568703172488
$ws.Range("K15").Value2 = 3
$ws.Range("M15").Value2 = 31.07813833333333
$ws.Range("N15").Value2 = 93.234415
$ws.Range("O15").Value2 = 0.2178407133195748
$ws.Range("P15").Value2 = 0.2178407133195749
$ws.Range("Q15").Value2 = 369.8864187378127
$ws.Range("R15").Value2 = 3328.977768640315
$ws.Range("S15").Value2 = 0.0362854451158113
$ws.Range("T15").Value2 = 0.03628544511581131

# Row 16
$ws.Range("E16").Value2 = 3
$ws.Range("G16").Value2 = 11.90182033333333
$ws.Range("H16").Value2 = 35.705461
$ws.Range("I16").Value2 = 0.166568703172488
$ws.Range("J16").Value2 = 0.166568703172488
$ws.Range("K16").Value2 = 3
$ws.Range("M16").Value2 = 91.56894199999999
$ws.Range("N16").Value2 = 274.706826
$ws.Range("O16").Value2 = 0.6418480872068143
$ws.Range("P16").Value2 = 0.6418480872068144
$ws.Range("Q16").Value2 = 1089.83709579742
$ws.Range("R16").Value2 = 9808.533862176786
$ws.Range("S16").Value2 = 0.106911803519781
$ws.Range("T16").Value2 = 0.1069118035197811

# Row 17
$ws.Range("E17").Value2 = 3
$ws.Range("G17").Value2 = 11.90182033333333
$ws.Range("H17").Value2 = 35.705461
$ws.Range("I17").Value2 = 0.166568703172488
$ws.Range("J17").Value2 = 0.166568703172488
$ws.Range("K17").Value2 = 3
$ws.Range("M17").Value2 = 8.822188333333335
$ws.Range("N17").Value2 = 26.466565
$ws.Range("O17").Value2 = 0.06183870407422939
$ws.Range("P17").Value2 = 0.06183870407422939
$ws.Range("Q17").Value2 = 105.0001004901628
$ws.Range("R17").Value2 = 945.0009044114651
$ws.Range("S17").Value2 = 0.01030039274351164
$ws.Range("T17").Value2 = 0.01030039274351164

